$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 36, shifting existing rows 36-88 down to 37-89.
$ws.Rows("36:36").Insert()

# Populate the newly inserted row 36 with its data (matches the style/pattern of
# the surrounding "Apio" records for Terminal Hortofrutícola Agro Chillán).
$ws.Range("A36").Value = 7
$ws.Range("B36").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C36").Value = "Ñuble"
$ws.Range("D36").Value = 44413
$ws.Range("E36").Value = 16
$ws.Range("F36").Value = 100112017
$ws.Range("G36").Value = "Apio"
$ws.Range("H36").Value = "Americana (o)"
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 120
$ws.Range("K36").Value = 8500
$ws.Range("L36").Value = 9000
$ws.Range("M36").Value = 8750
$ws.Range("N36").Value = "$/docena de matas"
$ws.Range("O36").Value = "Provincia del Elquí"
$ws.Range("P36").Value = 1458
$ws.Range("Q36").Value = 6
$ws.Range("R36").Value = "Hortaliza"
